# Fruta / hortaliza, semanal
# Insert two new weekly rows (new rows 259 and 260) into the "Limón" sheet,
# pushing the previously existing rows 259..282 down to 261..284.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right after the current row 258 (new rows 259 & 260).
$ws.Rows.Item(259).Insert()
$ws.Rows.Item(260).Insert()

# --- New row 259 ---
$ws.Cells.Item(259, 1).Value = 11
$ws.Cells.Item(259, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(259, 3).Value = "Bíobío"
$ws.Cells.Item(259, 4).Value = 44491
$ws.Cells.Item(259, 5).Value = 8
$ws.Cells.Item(259, 6).Value = "Fruta"
$ws.Cells.Item(259, 7).Value = 100102
$ws.Cells.Item(259, 8).Value = "Cítricos"
$ws.Cells.Item(259, 9).Value = 100102003
$ws.Cells.Item(259, 10).Value = "Limón"
$ws.Cells.Item(259, 11).Value = "Sin especificar"
$ws.Cells.Item(259, 12).Value = "1a amarillo"
$ws.Cells.Item(259, 13).Value = 180
$ws.Cells.Item(259, 14).Value = 5500
$ws.Cells.Item(259, 15).Value = 6000
$ws.Cells.Item(259, 16).Value = 5722
$ws.Cells.Item(259, 17).Value = "$/malla 16 kilos"
$ws.Cells.Item(259, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(259, 19).Value = 358
$ws.Cells.Item(259, 20).Value = 16

# --- New row 260 ---
$ws.Cells.Item(260, 1).Value = 11
$ws.Cells.Item(260, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(260, 3).Value = "Bíobío"
$ws.Cells.Item(260, 4).Value = 44491
$ws.Cells.Item(260, 5).Value = 8
$ws.Cells.Item(260, 6).Value = "Fruta"
$ws.Cells.Item(260, 7).Value = 100102
$ws.Cells.Item(260, 8).Value = "Cítricos"
$ws.Cells.Item(260, 9).Value = 100102003
$ws.Cells.Item(260, 10).Value = "Limón"
$ws.Cells.Item(260, 11).Value = "Sin especificar"
$ws.Cells.Item(260, 12).Value = "1a amarillo"
$ws.Cells.Item(260, 13).Value = 250
$ws.Cells.Item(260, 14).Value = 7000
$ws.Cells.Item(260, 15).Value = 7500
$ws.Cells.Item(260, 16).Value = 7300
$ws.Cells.Item(260, 17).Value = "$/malla 18 kilos"
$ws.Cells.Item(260, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(260, 19).Value = 406
$ws.Cells.Item(260, 20).Value = 18
